$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sub_materials_database")

$ws.Range("A1").Value = "Sub Material Name"
$ws.Range("B1").Value = "Chemical Composition"
